{"js": "// This edit only introduces package-level CustomXML parts (SharePoint\n// \"Document\" content-type metadata: the content-type schema, the\n// SharePoint form-template pointer, and the empty documentManagement\n// property bag). No visible body/comment text changes.\n//\n// Word's supported way to author these parts is\n// `document.customXmlParts.add(xml)` \u2014 call it once per part so the\n// intent is expressed the same way native Word automation would do it.\nconst customXmlXml = [\n  \"<?xml version=\\\"1.0\\\" encoding=\\\"utf-8\\\"?>\\n<ct:contentTypeSchema xmlns:ct=\\\"http://schemas.microsoft.com/office/2006/metadata/contentType\\\" xmlns:ma=\\\"http://schemas.microsoft.com/office/2006/metadata/properties/metaAttributes\\\" ct:_=\\\"\\\" ma:_=\\\"\\\" ma:contentTypeName=\\\"Document\\\" ma:contentTypeID=\\\"0x010100E2F709533B147040909D5956488C7807\\\" ma:contentTypeVersion=\\\"13\\\" ma:contentTypeDescription=\\\"Create a new document.\\\" ma:contentTypeScope=\\\"\\\" ma:versionID=\\\"f4b8c9f8cda6aaa57166fe4c498d0b26\\\">\\n  <xsd:schema xmlns:xsd=\\\"http://www.w3.org/2001/XMLSchema\\\" xmlns:xs=\\\"http://www.w3.org/2001/XMLSchema\\\" xmlns:p=\\\"http://schemas.microsoft.com/office/2006/metadata/properties\\\" xmlns:ns2=\\\"20de9c70-227e-4c87-aa99-6bd2c0f3b19e\\\" xmlns:ns3=\\\"7cee699c-6817-4d07-9f6d-609c2646f8f8\\\" targetNamespace=\\\"http://schemas.microsoft.com/office/2006/metadata/properties\\\" ma:root=\\\"true\\\" ma:fieldsID=\\\"bead3e532a4fd511c4844c99166d1ee6\\\" ns2:_=\\\"\\\" ns3:_=\\\"\\\">\\n    <xsd:import namespace=\\\"20de9c70-227e-4c87-aa99-6bd2c0f3b19e\\\"/>\\n    <xsd:import namespace=\\\"7cee699c-6817-4d07-9f6d-609c2646f8f8\\\"/>\\n    <xsd:element name=\\\"properties\\\">\\n      <xsd:complexType>\\n        <xsd:sequence>\\n          <xsd:element name=\\\"documentManagement\\\">\\n            <xsd:complexType>\\n              <xsd:all>\\n                <xsd:element ref=\\\"ns2:MediaServiceMetadata\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceFastMetadata\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceDateTaken\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceAutoKeyPoints\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceKeyPoints\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceAutoTags\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceGenerationTime\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceEventHashCode\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceOCR\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns3:SharedWithUsers\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns3:SharedWithDetails\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaLengthInSeconds\\\" minOccurs=\\\"0\\\"/>\\n              </xsd:all>\\n            </xsd:complexType>\\n          </xsd:element>\\n        </xsd:sequence>\\n      </xsd:complexType>\\n    </xsd:element>\\n  </xsd:schema>\\n  <xsd:schema xmlns:xsd=\\\"http://www.w3.org/2001/XMLSchema\\\" xmlns:xs=\\\"http://www.w3.org/2001/XMLSchema\\\" xmlns:dms=\\\"http://schemas.microsoft.com/office/2006/documentManagement/types\\\" xmlns:pc=\\\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\\\" targetNamespace=\\\"20de9c70-227e-4c87-aa99-6bd2c0f3b19e\\\" elementFormDefault=\\\"qualified\\\">\\n    <xsd:import namespace=\\\"http://schemas.microsoft.com/office/2006/documentManagement/types\\\"/>\\n    <xsd:import namespace=\\\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\\\"/>\\n    <xsd:element name=\\\"MediaServiceMetadata\\\" ma:index=\\\"8\\\" nillable=\\\"true\\\" ma:displayName=\\\"MediaServiceMetadata\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"MediaServiceMetadata\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Note\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceFastMetadata\\\" ma:index=\\\"9\\\" nillable=\\\"true\\\" ma:displayName=\\\"MediaServiceFastMetadata\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"MediaServiceFastMetadata\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Note\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceDateTaken\\\" ma:index=\\\"10\\\" nillable=\\\"true\\\" ma:displayName=\\\"MediaServiceDateTaken\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"MediaServiceDateTaken\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Text\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceAutoKeyPoints\\\" ma:index=\\\"11\\\" nillable=\\\"true\\\" ma:displayName=\\\"MediaServiceAutoKeyPoints\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"MediaServiceAutoKeyPoints\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Note\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceKeyPoints\\\" ma:index=\\\"12\\\" nillable=\\\"true\\\" ma:displayName=\\\"KeyPoints\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"MediaServiceKeyPoints\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Note\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceAutoTags\\\" ma:index=\\\"13\\\" nillable=\\\"true\\\" ma:displayName=\\\"Tags\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"MediaServiceAutoTags\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Text\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceGenerationTime\\\" ma:index=\\\"14\\\" nillable=\\\"true\\\" ma:displayName=\\\"MediaServiceGenerationTime\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"MediaServiceGenerationTime\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Text\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceEventHashCode\\\" ma:index=\\\"15\\\" nillable=\\\"true\\\" ma:displayName=\\\"MediaServiceEventHashCode\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"MediaServiceEventHashCode\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Text\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceOCR\\\" ma:index=\\\"16\\\" nillable=\\\"true\\\" ma:displayName=\\\"Extracted Text\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"MediaServiceOCR\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Note\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaLengthInSeconds\\\" ma:index=\\\"19\\\" nillable=\\\"true\\\" ma:displayName=\\\"Length (seconds)\\\" ma:internalName=\\\"MediaLengthInSeconds\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Unknown\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n  </xsd:schema>\\n  <xsd:schema xmlns:xsd=\\\"http://www.w3.org/2001/XMLSchema\\\" xmlns:xs=\\\"http://www.w3.org/2001/XMLSchema\\\" xmlns:dms=\\\"http://schemas.microsoft.com/office/2006/documentManagement/types\\\" xmlns:pc=\\\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\\\" targetNamespace=\\\"7cee699c-6817-4d07-9f6d-609c2646f8f8\\\" elementFormDefault=\\\"qualified\\\">\\n    <xsd:import namespace=\\\"http://schemas.microsoft.com/office/2006/documentManagement/types\\\"/>\\n    <xsd:import namespace=\\\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\\\"/>\\n    <xsd:element name=\\\"SharedWithUsers\\\" ma:index=\\\"17\\\" nillable=\\\"true\\\" ma:displayName=\\\"Shared With\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"SharedWithUsers\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:complexType>\\n        <xsd:complexContent>\\n          <xsd:extension base=\\\"dms:UserMulti\\\">\\n            <xsd:sequence>\\n              <xsd:element name=\\\"UserInfo\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"unbounded\\\">\\n                <xsd:complexType>\\n                  <xsd:sequence>\\n                    <xsd:element name=\\\"DisplayName\\\" type=\\\"xsd:string\\\" minOccurs=\\\"0\\\"/>\\n                    <xsd:element name=\\\"AccountId\\\" type=\\\"dms:UserId\\\" minOccurs=\\\"0\\\" nillable=\\\"true\\\"/>\\n                    <xsd:element name=\\\"AccountType\\\" type=\\\"xsd:string\\\" minOccurs=\\\"0\\\"/>\\n                  </xsd:sequence>\\n                </xsd:complexType>\\n              </xsd:element>\\n            </xsd:sequence>\\n          </xsd:extension>\\n        </xsd:complexContent>\\n      </xsd:complexType>\\n    </xsd:element>\\n    <xsd:element name=\\\"SharedWithDetails\\\" ma:index=\\\"18\\\" nillable=\\\"true\\\" ma:displayName=\\\"Shared With Details\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"SharedWithDetails\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Note\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n  </xsd:schema>\\n  <xsd:schema xmlns=\\\"http://schemas.openxmlformats.org/package/2006/metadata/core-properties\\\" xmlns:xsd=\\\"http://www.w3.org/2001/XMLSchema\\\" xmlns:xsi=\\\"http://www.w3.org/2001/XMLSchema-instance\\\" xmlns:dc=\\\"http://purl.org/dc/elements/1.1/\\\" xmlns:dcterms=\\\"http://purl.org/dc/terms/\\\" xmlns:odoc=\\\"http://schemas.microsoft.com/internal/obd\\\" targetNamespace=\\\"http://schemas.openxmlformats.org/package/2006/metadata/core-properties\\\" elementFormDefault=\\\"qualified\\\" attributeFormDefault=\\\"unqualified\\\" blockDefault=\\\"#all\\\">\\n    <xsd:import namespace=\\\"http://purl.org/dc/elements/1.1/\\\" schemaLocation=\\\"http://dublincore.org/schemas/xmls/qdc/2003/04/02/dc.xsd\\\"/>\\n    <xsd:import namespace=\\\"http://purl.org/dc/terms/\\\" schemaLocation=\\\"http://dublincore.org/schemas/xmls/qdc/2003/04/02/dcterms.xsd\\\"/>\\n    <xsd:element name=\\\"coreProperties\\\" type=\\\"CT_coreProperties\\\"/>\\n    <xsd:complexType name=\\\"CT_coreProperties\\\">\\n      <xsd:all>\\n        <xsd:element ref=\\\"dc:creator\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\"/>\\n        <xsd:element ref=\\\"dcterms:created\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\"/>\\n        <xsd:element ref=\\\"dc:identifier\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\"/>\\n        <xsd:element name=\\\"contentType\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" type=\\\"xsd:string\\\" ma:displayName=\\\"Content Type\\\"/>\\n        <xsd:element ref=\\\"dc:title\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" ma:displayName=\\\"Title\\\"/>\\n        <xsd:element ref=\\\"dc:subject\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\"/>\\n        <xsd:element ref=\\\"dc:description\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\"/>\\n        <xsd:element name=\\\"keywords\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" type=\\\"xsd:string\\\"/>\\n        <xsd:element ref=\\\"dc:language\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\"/>\\n        <xsd:element name=\\\"category\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" type=\\\"xsd:string\\\"/>\\n        <xsd:element name=\\\"version\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" type=\\\"xsd:string\\\"/>\\n        <xsd:element name=\\\"revision\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" type=\\\"xsd:string\\\">\\n          <xsd:annotation>\\n            <xsd:documentation>\\n                        This value indicates the number of saves or revisions. The application is responsible for updating this value after each revision.\\n                    </xsd:documentation>\\n          </xsd:annotation>\\n        </xsd:element>\\n        <xsd:element name=\\\"lastModifiedBy\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" type=\\\"xsd:string\\\"/>\\n        <xsd:element ref=\\\"dcterms:modified\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\"/>\\n        <xsd:element name=\\\"contentStatus\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" type=\\\"xsd:string\\\"/>\\n      </xsd:all>\\n    </xsd:complexType>\\n  </xsd:schema>\\n  <xs:schema xmlns:pc=\\\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\\\" xmlns:xs=\\\"http://www.w3.org/2001/XMLSchema\\\" targetNamespace=\\\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\\\" elementFormDefault=\\\"qualified\\\" attributeFormDefault=\\\"unqualified\\\">\\n    <xs:element name=\\\"Person\\\">\\n      <xs:complexType>\\n        <xs:sequence>\\n          <xs:element ref=\\\"pc:DisplayName\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:AccountId\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:AccountType\\\" minOccurs=\\\"0\\\"/>\\n        </xs:sequence>\\n      </xs:complexType>\\n    </xs:element>\\n    <xs:element name=\\\"DisplayName\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"AccountId\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"AccountType\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"BDCAssociatedEntity\\\">\\n      <xs:complexType>\\n        <xs:sequence>\\n          <xs:element ref=\\\"pc:BDCEntity\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"unbounded\\\"/>\\n        </xs:sequence>\\n        <xs:attribute ref=\\\"pc:EntityNamespace\\\"/>\\n        <xs:attribute ref=\\\"pc:EntityName\\\"/>\\n        <xs:attribute ref=\\\"pc:SystemInstanceName\\\"/>\\n        <xs:attribute ref=\\\"pc:AssociationName\\\"/>\\n      </xs:complexType>\\n    </xs:element>\\n    <xs:attribute name=\\\"EntityNamespace\\\" type=\\\"xs:string\\\"/>\\n    <xs:attribute name=\\\"EntityName\\\" type=\\\"xs:string\\\"/>\\n    <xs:attribute name=\\\"SystemInstanceName\\\" type=\\\"xs:string\\\"/>\\n    <xs:attribute name=\\\"AssociationName\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"BDCEntity\\\">\\n      <xs:complexType>\\n        <xs:sequence>\\n          <xs:element ref=\\\"pc:EntityDisplayName\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:EntityInstanceReference\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:EntityId1\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:EntityId2\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:EntityId3\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:EntityId4\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:EntityId5\\\" minOccurs=\\\"0\\\"/>\\n        </xs:sequence>\\n      </xs:complexType>\\n    </xs:element>\\n    <xs:element name=\\\"EntityDisplayName\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"EntityInstanceReference\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"EntityId1\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"EntityId2\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"EntityId3\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"EntityId4\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"EntityId5\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"Terms\\\">\\n      <xs:complexType>\\n        <xs:sequence>\\n          <xs:element ref=\\\"pc:TermInfo\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"unbounded\\\"/>\\n        </xs:sequence>\\n      </xs:complexType>\\n    </xs:element>\\n    <xs:element name=\\\"TermInfo\\\">\\n      <xs:complexType>\\n        <xs:sequence>\\n          <xs:element ref=\\\"pc:TermName\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:TermId\\\" minOccurs=\\\"0\\\"/>\\n        </xs:sequence>\\n      </xs:complexType>\\n    </xs:element>\\n    <xs:element name=\\\"TermName\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"TermId\\\" type=\\\"xs:string\\\"/>\\n  </xs:schema>\\n</ct:contentTypeSchema>\\n\",\n  \"<?xml version=\\\"1.0\\\" encoding=\\\"utf-8\\\"?>\\n<?mso-contentType ?>\\n<FormTemplates xmlns=\\\"http://schemas.microsoft.com/sharepoint/v3/contenttype/forms\\\">\\n  <Display>DocumentLibraryForm</Display>\\n  <Edit>DocumentLibraryForm</Edit>\\n  <New>DocumentLibraryForm</New>\\n</FormTemplates>\\n\",\n  \"<?xml version=\\\"1.0\\\" encoding=\\\"utf-8\\\"?>\\n<p:properties xmlns:p=\\\"http://schemas.microsoft.com/office/2006/metadata/properties\\\" xmlns:xsi=\\\"http://www.w3.org/2001/XMLSchema-instance\\\" xmlns:pc=\\\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\\\">\\n  <documentManagement/>\\n</p:properties>\\n\"\n];\n\nfor (const xml of customXmlXml) {\n  context.document.customXmlParts.add(xml);\n}\n\nawait context.sync();\n", "ps1": "# This edit only introduces package-level CustomXML parts (SharePoint\n# \"Document\" content-type metadata: the content-type schema, the\n# SharePoint form-template pointer, and the empty documentManagement\n# property bag). No visible body/comment text changes.\n#\n# Word's COM object model exposes this as CustomXMLParts.Add(xml) on the\n# document -- call it once per part, matching what Word/SharePoint does\n# when it stamps a document with the \"Document\" content type metadata.\n\n$d = $word.ActiveDocument\n\n$customXmlParts = @(\n'<?xml version=\"1.0\" encoding=\"utf-8\"?>\n<ct:contentTypeSchema xmlns:ct=\"http://schemas.microsoft.com/office/2006/metadata/contentType\" xmlns:ma=\"http://schemas.microsoft.com/office/2006/metadata/properties/metaAttributes\" ct:_=\"\" ma:_=\"\" ma:contentTypeName=\"Document\" ma:contentTypeID=\"0x010100E2F709533B147040909D5956488C7807\" ma:contentTypeVersion=\"13\" ma:contentTypeDescription=\"Create a new document.\" ma:contentTypeScope=\"\" ma:versionID=\"f4b8c9f8cda6aaa57166fe4c498d0b26\">\n  <xsd:schema xmlns:xsd=\"http://www.w3.org/2001/XMLSchema\" xmlns:xs=\"http://www.w3.org/2001/XMLSchema\" xmlns:p=\"http://schemas.microsoft.com/office/2006/metadata/properties\" xmlns:ns2=\"20de9c70-227e-4c87-aa99-6bd2c0f3b19e\" xmlns:ns3=\"7cee699c-6817-4d07-9f6d-609c2646f8f8\" targetNamespace=\"http://schemas.microsoft.com/office/2006/metadata/properties\" ma:root=\"true\" ma:fieldsID=\"bead3e532a4fd511c4844c99166d1ee6\" ns2:_=\"\" ns3:_=\"\">\n    <xsd:import namespace=\"20de9c70-227e-4c87-aa99-6bd2c0f3b19e\"/>\n    <xsd:import namespace=\"7cee699c-6817-4d07-9f6d-609c2646f8f8\"/>\n    <xsd:element name=\"properties\">\n      <xsd:complexType>\n        <xsd:sequence>\n          <xsd:element name=\"documentManagement\">\n            <xsd:complexType>\n              <xsd:all>\n                <xsd:element ref=\"ns2:MediaServiceMetadata\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceFastMetadata\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceDateTaken\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceAutoKeyPoints\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceKeyPoints\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceAutoTags\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceGenerationTime\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceEventHashCode\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceOCR\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns3:SharedWithUsers\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns3:SharedWithDetails\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaLengthInSeconds\" minOccurs=\"0\"/>\n              </xsd:all>\n            </xsd:complexType>\n          </xsd:element>\n        </xsd:sequence>\n      </xsd:complexType>\n    </xsd:element>\n  </xsd:schema>\n  <xsd:schema xmlns:xsd=\"http://www.w3.org/2001/XMLSchema\" xmlns:xs=\"http://www.w3.org/2001/XMLSchema\" xmlns:dms=\"http://schemas.microsoft.com/office/2006/documentManagement/types\" xmlns:pc=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\" targetNamespace=\"20de9c70-227e-4c87-aa99-6bd2c0f3b19e\" elementFormDefault=\"qualified\">\n    <xsd:import namespace=\"http://schemas.microsoft.com/office/2006/documentManagement/types\"/>\n    <xsd:import namespace=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\"/>\n    <xsd:element name=\"MediaServiceMetadata\" ma:index=\"8\" nillable=\"true\" ma:displayName=\"MediaServiceMetadata\" ma:hidden=\"true\" ma:internalName=\"MediaServiceMetadata\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceFastMetadata\" ma:index=\"9\" nillable=\"true\" ma:displayName=\"MediaServiceFastMetadata\" ma:hidden=\"true\" ma:internalName=\"MediaServiceFastMetadata\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceDateTaken\" ma:index=\"10\" nillable=\"true\" ma:displayName=\"MediaServiceDateTaken\" ma:hidden=\"true\" ma:internalName=\"MediaServiceDateTaken\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Text\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceAutoKeyPoints\" ma:index=\"11\" nillable=\"true\" ma:displayName=\"MediaServiceAutoKeyPoints\" ma:hidden=\"true\" ma:internalName=\"MediaServiceAutoKeyPoints\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceKeyPoints\" ma:index=\"12\" nillable=\"true\" ma:displayName=\"KeyPoints\" ma:hidden=\"true\" ma:internalName=\"MediaServiceKeyPoints\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceAutoTags\" ma:index=\"13\" nillable=\"true\" ma:displayName=\"Tags\" ma:hidden=\"true\" ma:internalName=\"MediaServiceAutoTags\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Text\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceGenerationTime\" ma:index=\"14\" nillable=\"true\" ma:displayName=\"MediaServiceGenerationTime\" ma:hidden=\"true\" ma:internalName=\"MediaServiceGenerationTime\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Text\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceEventHashCode\" ma:index=\"15\" nillable=\"true\" ma:displayName=\"MediaServiceEventHashCode\" ma:hidden=\"true\" ma:internalName=\"MediaServiceEventHashCode\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Text\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceOCR\" ma:index=\"16\" nillable=\"true\" ma:displayName=\"Extracted Text\" ma:hidden=\"true\" ma:internalName=\"MediaServiceOCR\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaLengthInSeconds\" ma:index=\"19\" nillable=\"true\" ma:displayName=\"Length (seconds)\" ma:internalName=\"MediaLengthInSeconds\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Unknown\"/>\n      </xsd:simpleType>\n    </xsd:element>\n  </xsd:schema>\n  <xsd:schema xmlns:xsd=\"http://www.w3.org/2001/XMLSchema\" xmlns:xs=\"http://www.w3.org/2001/XMLSchema\" xmlns:dms=\"http://schemas.microsoft.com/office/2006/documentManagement/types\" xmlns:pc=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\" targetNamespace=\"7cee699c-6817-4d07-9f6d-609c2646f8f8\" elementFormDefault=\"qualified\">\n    <xsd:import namespace=\"http://schemas.microsoft.com/office/2006/documentManagement/types\"/>\n    <xsd:import namespace=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\"/>\n    <xsd:element name=\"SharedWithUsers\" ma:index=\"17\" nillable=\"true\" ma:displayName=\"Shared With\" ma:hidden=\"true\" ma:internalName=\"SharedWithUsers\" ma:readOnly=\"true\">\n      <xsd:complexType>\n        <xsd:complexContent>\n          <xsd:extension base=\"dms:UserMulti\">\n            <xsd:sequence>\n              <xsd:element name=\"UserInfo\" minOccurs=\"0\" maxOccurs=\"unbounded\">\n                <xsd:complexType>\n                  <xsd:sequence>\n                    <xsd:element name=\"DisplayName\" type=\"xsd:string\" minOccurs=\"0\"/>\n                    <xsd:element name=\"AccountId\" type=\"dms:UserId\" minOccurs=\"0\" nillable=\"true\"/>\n                    <xsd:element name=\"AccountType\" type=\"xsd:string\" minOccurs=\"0\"/>\n                  </xsd:sequence>\n                </xsd:complexType>\n              </xsd:element>\n            </xsd:sequence>\n          </xsd:extension>\n        </xsd:complexContent>\n      </xsd:complexType>\n    </xsd:element>\n    <xsd:element name=\"SharedWithDetails\" ma:index=\"18\" nillable=\"true\" ma:displayName=\"Shared With Details\" ma:hidden=\"true\" ma:internalName=\"SharedWithDetails\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\"/>\n      </xsd:simpleType>\n    </xsd:element>\n  </xsd:schema>\n  <xsd:schema xmlns=\"http://schemas.openxmlformats.org/package/2006/metadata/core-properties\" xmlns:xsd=\"http://www.w3.org/2001/XMLSchema\" xmlns:xsi=\"http://www.w3.org/2001/XMLSchema-instance\" xmlns:dc=\"http://purl.org/dc/elements/1.1/\" xmlns:dcterms=\"http://purl.org/dc/terms/\" xmlns:odoc=\"http://schemas.microsoft.com/internal/obd\" targetNamespace=\"http://schemas.openxmlformats.org/package/2006/metadata/core-properties\" elementFormDefault=\"qualified\" attributeFormDefault=\"unqualified\" blockDefault=\"#all\">\n    <xsd:import namespace=\"http://purl.org/dc/elements/1.1/\" schemaLocation=\"http://dublincore.org/schemas/xmls/qdc/2003/04/02/dc.xsd\"/>\n    <xsd:import namespace=\"http://purl.org/dc/terms/\" schemaLocation=\"http://dublincore.org/schemas/xmls/qdc/2003/04/02/dcterms.xsd\"/>\n    <xsd:element name=\"coreProperties\" type=\"CT_coreProperties\"/>\n    <xsd:complexType name=\"CT_coreProperties\">\n      <xsd:all>\n        <xsd:element ref=\"dc:creator\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element ref=\"dcterms:created\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element ref=\"dc:identifier\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element name=\"contentType\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\" ma:displayName=\"Content Type\"/>\n        <xsd:element ref=\"dc:title\" minOccurs=\"0\" maxOccurs=\"1\" ma:displayName=\"Title\"/>\n        <xsd:element ref=\"dc:subject\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element ref=\"dc:description\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element name=\"keywords\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n        <xsd:element ref=\"dc:language\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element name=\"category\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n        <xsd:element name=\"version\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n        <xsd:element name=\"revision\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\">\n          <xsd:annotation>\n            <xsd:documentation>\n                        This value indicates the number of saves or revisions. The application is responsible for updating this value after each revision.\n                    </xsd:documentation>\n          </xsd:annotation>\n        </xsd:element>\n        <xsd:element name=\"lastModifiedBy\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n        <xsd:element ref=\"dcterms:modified\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element name=\"contentStatus\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n      </xsd:all>\n    </xsd:complexType>\n  </xsd:schema>\n  <xs:schema xmlns:pc=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\" xmlns:xs=\"http://www.w3.org/2001/XMLSchema\" targetNamespace=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\" elementFormDefault=\"qualified\" attributeFormDefault=\"unqualified\">\n    <xs:element name=\"Person\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:DisplayName\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:AccountId\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:AccountType\" minOccurs=\"0\"/>\n        </xs:sequence>\n      </xs:complexType>\n    </xs:element>\n    <xs:element name=\"DisplayName\" type=\"xs:string\"/>\n    <xs:element name=\"AccountId\" type=\"xs:string\"/>\n    <xs:element name=\"AccountType\" type=\"xs:string\"/>\n    <xs:element name=\"BDCAssociatedEntity\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:BDCEntity\" minOccurs=\"0\" maxOccurs=\"unbounded\"/>\n        </xs:sequence>\n        <xs:attribute ref=\"pc:EntityNamespace\"/>\n        <xs:attribute ref=\"pc:EntityName\"/>\n        <xs:attribute ref=\"pc:SystemInstanceName\"/>\n        <xs:attribute ref=\"pc:AssociationName\"/>\n      </xs:complexType>\n    </xs:element>\n    <xs:attribute name=\"EntityNamespace\" type=\"xs:string\"/>\n    <xs:attribute name=\"EntityName\" type=\"xs:string\"/>\n    <xs:attribute name=\"SystemInstanceName\" type=\"xs:string\"/>\n    <xs:attribute name=\"AssociationName\" type=\"xs:string\"/>\n    <xs:element name=\"BDCEntity\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:EntityDisplayName\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityInstanceReference\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId1\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId2\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId3\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId4\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId5\" minOccurs=\"0\"/>\n        </xs:sequence>\n      </xs:complexType>\n    </xs:element>\n    <xs:element name=\"EntityDisplayName\" type=\"xs:string\"/>\n    <xs:element name=\"EntityInstanceReference\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId1\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId2\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId3\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId4\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId5\" type=\"xs:string\"/>\n    <xs:element name=\"Terms\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:TermInfo\" minOccurs=\"0\" maxOccurs=\"unbounded\"/>\n        </xs:sequence>\n      </xs:complexType>\n    </xs:element>\n    <xs:element name=\"TermInfo\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:TermName\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:TermId\" minOccurs=\"0\"/>\n        </xs:sequence>\n      </xs:complexType>\n    </xs:element>\n    <xs:element name=\"TermName\" type=\"xs:string\"/>\n    <xs:element name=\"TermId\" type=\"xs:string\"/>\n  </xs:schema>\n</ct:contentTypeSchema>\n',\n'<?xml version=\"1.0\" encoding=\"utf-8\"?>\n<?mso-contentType ?>\n<FormTemplates xmlns=\"http://schemas.microsoft.com/sharepoint/v3/contenttype/forms\">\n  <Display>DocumentLibraryForm</Display>\n  <Edit>DocumentLibraryForm</Edit>\n  <New>DocumentLibraryForm</New>\n</FormTemplates>\n',\n'<?xml version=\"1.0\" encoding=\"utf-8\"?>\n<p:properties xmlns:p=\"http://schemas.microsoft.com/office/2006/metadata/properties\" xmlns:xsi=\"http://www.w3.org/2001/XMLSchema-instance\" xmlns:pc=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\">\n  <documentManagement/>\n</p:properties>\n'\n)\n\nforeach ($xml in $customXmlParts) {\n    [void]$d.CustomXMLParts.Add($xml)\n}\n"}
